$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Rush Hour tendrá 6 clases..." -> split "Hour" into its own run so the
#    sentence reads as three runs: "Rush " | "Hour" | " tendrá 6 clases...".
#    Toggling a character property on the narrow sub-range forces Word to
#    break the parent run at those exact boundaries while leaving the
#    run-level formatting (rFonts/lang) untouched once the toggle is undone.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Hour", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Bold = 1
$rng.Bold = 0

# ---------------------------------------------------------------------------
# 2) "...afectación que ocurra en estas 2 se vera reflejada..." -> split
#    "vera" into its own run the same way: "...se " | "vera" | " reflejada...".
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("vera", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Bold = 1
$rng.Bold = 0

# ---------------------------------------------------------------------------
# 3) "Grandma" bullet: merge the trailing ":" run and the " Esta clase..."
#    run back into a single run reading ": Esta clase...". A same-text
#    Find/Replace across the run boundary collapses them into one run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$grandmaText = ": Esta clase corresponde a un NPC especial, es por ello que hereda de la clase descrita anteriormente, la única diferencia es que las abuelitas pueden solicitarle un puesto tanto a los NPC como al jugador."
$rng.Find.Execute($grandmaText, $true, $false, $false, $false, $false, $true, 1, $false, $grandmaText, 2)

# ---------------------------------------------------------------------------
# 4) Drop the trailing empty paragraphs, the "El motivo..." paragraph, and
#    the final empty paragraph that followed the class-diagram picture,
#    leaving the picture's paragraph as the last one in the body.
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$total = $paras.Count
$drawingParaIndex = 0
for ($i = 1; $i -le $total; $i++) {
    if ($paras.Item($i).Range.InlineShapes.Count -gt 0) {
        $drawingParaIndex = $i
    }
}
if ($drawingParaIndex -gt 0 -and $drawingParaIndex -lt $total) {
    $startPara = $paras.Item($drawingParaIndex + 1)
    $endPara = $paras.Item($total)
    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete()
}
